# Update "想去人数" (F column) values on the "展览" sheet and the
# aggregated "全部类型" sheet, per the commit's regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (row -> new F value) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 98
$ws1.Range("F4").Value  = 407
$ws1.Range("F6").Value  = 129
$ws1.Range("F7").Value  = 1083
$ws1.Range("F8").Value  = 363
$ws1.Range("F9").Value  = 185
$ws1.Range("F13").Value = 366
$ws1.Range("F19").Value = 987
$ws1.Range("F20").Value = 445
$ws1.Range("F22").Value = 80
$ws1.Range("F23").Value = 372

# --- Sheet "全部类型" (same events, different row numbers) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 98
$ws4.Range("F6").Value  = 407
$ws4.Range("F8").Value  = 129
$ws4.Range("F9").Value  = 1083
$ws4.Range("F10").Value = 363
$ws4.Range("F11").Value = 185
$ws4.Range("F20").Value = 366
$ws4.Range("F26").Value = 987
$ws4.Range("F27").Value = 445
$ws4.Range("F31").Value = 80
$ws4.Range("F32").Value = 372
